$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.949.99'
$ws.Range('E2').Value = '  -1.93%  '
$ws.Range('D3').Value = '2.317.16'
$ws.Range('E3').Value = '  -4.19%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'549.16"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').Value = "'131.66"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.76%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -2.60%  '
$ws.Range('D9').Value = '2.316.23'
$ws.Range('E9').Value = '  -4.10%  '
$ws.Range('E10').Value = '  -3.05%  '
$ws.Range('D11').Value = "'5.57"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.04%  '
$ws.Range('E13').Value = '  -5.15%  '
$ws.Range('D14').Value = "'24.04"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.00%  '
$ws.Range('D15').Value = '2.731.49'
$ws.Range('E15').Value = '  -4.18%  '
$ws.Range('D16').Value = '58.851.19'
$ws.Range('E16').Value = '  -1.95%  '
$ws.Range('E17').Value = '  -2.89%  '
$ws.Range('D18').Value = '2.335.69'
$ws.Range('E18').Value = '  -4.36%  '
$ws.Range('E19').Value = '  -4.48%  '
$ws.Range('D20').Value = "'4.33"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.66%  '
$ws.Range('D21').Value = "'316.26"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.39%  '
$ws.Range('D22').Value = "'6.51"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.14%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = "'63.59"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.91%  '
$ws.Range('E25').Value = '  -3.97%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').Value = "'8.09"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.33%  '
$ws.Range('D28').Value = "'1.33"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.61%  '
$ws.Range('D29').Value = "'1.76"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').Value = "'169.57"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('E31').Value = '  -5.37%  '
$ws.Range('E32').Value = '  +4.61%  '
$ws.Range('D33').Value = "'5.81"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.81%  '
$ws.Range('E34').Value = '  -4.27%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').Value = '  -3.73%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').Value = "'1.26"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.06%  '
$ws.Range('E39').Value = '  -5.45%  '
$ws.Range('D40').Value = "'38.21"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('E41').Value = '  -4.56%  '
$ws.Range('D42').Value = "'303.01"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.93%  '
$ws.Range('D43').Value = "'141.73"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.78%  '
$ws.Range('E44').Value = '  -5.49%  '
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = "'18.74"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.36%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = "'0.561"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.07%  '
$ws.Range('E49').Value = '  -2.75%  '
$ws.Range('D50').Value = "'16.71"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.25%  '
$ws.Range('D51').Value = "'11.02"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.20%  '
